$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text that looks numeric (e.g. "577.12").
# Excel Value auto-coerces those to real numbers, changing the stored cell
# type away from text. Force text storage via NumberFormat, then reset the
# style back to Normal so no stray style index lingers on the cell.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "67.839.19"
$ws.Range("E2").Value = "  -6.51%  "
# Row 3
Set-TextValue $ws.Range("D3") "3.694.71"
$ws.Range("E3").Value = "  -5.96%  "
# Row 4
$ws.Range("E4").Value = "  +0.12%  "
# Row 5
Set-TextValue $ws.Range("D5") "577.12"
$ws.Range("E5").Value = "  -3.82%  "
# Row 6
Set-TextValue $ws.Range("D6") "175.03"
$ws.Range("E6").Value = "  +2.85%  "
# Row 7
Set-TextValue $ws.Range("D7") "3.686.83"
$ws.Range("E7").Value = "  -6.05%  "
# Row 8
Set-TextValue $ws.Range("D8") "0.623"
$ws.Range("E8").Value = "  -8.54%  "
# Row 9
Set-TextValue $ws.Range("D9") "1.00"
$ws.Range("E9").Value = "  +0.20%  "
# Row 10
$ws.Range("E10").Value = "  -9.85%  "
# Row 11
$ws.Range("E11").Value = "  -13.52%  "
# Row 12
Set-TextValue $ws.Range("D12") "51.33"
$ws.Range("E12").Value = "  -7.49%  "
# Row 14
Set-TextValue $ws.Range("D14") "10.37"
$ws.Range("E14").Value = "  -9.54%  "
# Row 15
Set-TextValue $ws.Range("D15") "4.286.11"
$ws.Range("E15").Value = "  -5.95%  "
# Row 16
Set-TextValue $ws.Range("D16") "3.689.79"
$ws.Range("E16").Value = "  -6.15%  "
# Row 17
Set-TextValue $ws.Range("D17") "19.26"
$ws.Range("E17").Value = "  -9.62%  "
# Row 18
$ws.Range("E18").Value = "  -3.20%  "
# Row 19
Set-TextValue $ws.Range("D19") "12.75"
$ws.Range("E19").Value = "  -9.67%  "
# Row 20
$ws.Range("E20").Value = "  -9.36%  "
# Row 21
Set-TextValue $ws.Range("D21") "67.588.18"
$ws.Range("E21").Value = "  -6.74%  "
# Row 22
Set-TextValue $ws.Range("D22") "404.10"
$ws.Range("E22").Value = "  -9.57%  "
# Row 23
Set-TextValue $ws.Range("D23") "4.44"
$ws.Range("E23").Value = "  -7.36%  "
# Row 24
Set-TextValue $ws.Range("D24") "87.84"
$ws.Range("E24").Value = "  -7.92%  "
# Row 25
$ws.Range("E25").Value = "  -8.83%  "
# Row 26
Set-TextValue $ws.Range("D26") "12.64"
$ws.Range("E26").Value = "  -9.75%  "
# Row 27
Set-TextValue $ws.Range("D27") "10.70"
$ws.Range("E27").Value = "  -3.46%  "
# Row 28
$ws.Range("E28").Value = "  +1.25%  "
# Row 29
Set-TextValue $ws.Range("D29") "3.79"
$ws.Range("E29").Value = "  -11.00%  "
# Row 30
Set-TextValue $ws.Range("D30") "9.43"
$ws.Range("E30").Value = "  -8.42%  "
# Row 31
Set-TextValue $ws.Range("D31") "32.38"
$ws.Range("E31").Value = "  -9.26%  "
# Row 32
Set-TextValue $ws.Range("D32") "7.38"
$ws.Range("E32").Value = "  -5.95%  "
# Row 33
Set-TextValue $ws.Range("D33") "12.36"
$ws.Range("E33").Value = "  -10.37%  "
# Row 34
Set-TextValue $ws.Range("D34") "610.07"
$ws.Range("E34").Value = "  -2.23%  "
# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D35") "0.115"
$ws.Range("E35").Value = "  -9.34%  "
# Row 36
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D36") "64.65"
$ws.Range("E36").Value = "  -5.94%  "
# Row 37
Set-TextValue $ws.Range("D37") "42.74"
$ws.Range("E37").Value = "  -15.57%  "
# Row 38
Set-TextValue $ws.Range("D38") "0.0₃0879"
$ws.Range("E38").Value = "  -10.47%  "
# Row 39
$ws.Range("E39").Value = "  +0.25%  "
# Row 40
Set-TextValue $ws.Range("D40") "0.393"
$ws.Range("E40").Value = "  -7.51%  "
# Row 41
$ws.Range("E41").Value = "  +0.06%  "
# Row 42
Set-TextValue $ws.Range("D42") "0.134"
$ws.Range("E42").Value = "  -7.05%  "
# Row 43
Set-TextValue $ws.Range("D43") "2.74"
$ws.Range("E43").Value = "  +6.55%  "
# Row 44
$ws.Range("E44").Value = "  -10.59%  "
# Row 45
Set-TextValue $ws.Range("D45") "0.0431"
$ws.Range("E45").Value = "  -9.40%  "
# Row 46
Set-TextValue $ws.Range("D46") "2.84"
$ws.Range("E46").Value = "  -12.17%  "
# Row 47
Set-TextValue $ws.Range("D47") "9.16"
$ws.Range("E47").Value = "  -12.49%  "
# Row 48
Set-TextValue $ws.Range("D48") "2.788.88"
$ws.Range("E48").Value = "  -1.38%  "
# Row 49
$ws.Range("E49").Value = "  -9.62%  "
# Row 50
Set-TextValue $ws.Range("D50") "2.68"
$ws.Range("E50").Value = "  -5.99%  "
# Row 51
Set-TextValue $ws.Range("D51") "3.02"
$ws.Range("E51").Value = "  -10.18%  "
